# Swap the roster data for row 9 ("JaVale McGee") and row 10 ("Dāvis Bertāns")
# so that Dāvis Bertāns's stats end up on row 9 and JaVale McGee's stats end
# up on row 10 (i.e. the two players trade places in the roster table).
#
# We do this with Copy / PasteSpecial (values) instead of literally retyping
# every value, so that cell data types (e.g. numeric-looking text such as
# "6" / "14" in the Exp column) and formatting are carried over faithfully
# instead of being reinterpreted (typing "6" directly would turn it into a
# real number instead of text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("B100:K100")
$row9 = $ws.Range("B9:K9")
$row10 = $ws.Range("B10:K10")

# Stash row 9's current contents (JaVale McGee) in a scratch area.
$row9.Copy()
$scratch.PasteSpecial(-4163)

# Move row 10's contents (Dāvis Bertāns) into row 9.
# Clear first so that cells which are blank in row 10 (e.g. the College
# column) actually become blank in row 9 rather than keeping stale data.
$row9.Clear()
$row10.Copy()
$row9.PasteSpecial(-4163)

# Move the stashed original row 9 contents (JaVale McGee) into row 10.
$row10.Clear()
$scratch.Copy()
$row10.PasteSpecial(-4163)

# Remove the scratch data/formatting so it doesn't linger in the sheet.
$scratch.Clear()

# Clearing the ranges above also wiped the hyperlink cell style from K9/K10;
# restore it so the bbref url cells keep their usual hyperlink formatting.
$ws.Range("K9").Style = "Hyperlink"
$ws.Range("K10").Style = "Hyperlink"
